# Apply the commit: insert 3 new data rows (weekly price records) right before
# the former row 817, shifting the remaining rows of the "Cebolla" sub-table
# down by three rows (old A1:R858 -> new A1:R861).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows at 817, pushing existing rows 817-858 down to 820-861.
$ws.Rows("817:819").Insert()

# --- New row 817 ---
$ws.Cells.Item(817, 1).Value2 = 5
$ws.Cells.Item(817, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(817, 3).Value2 = "Maule"
$ws.Cells.Item(817, 4).Value2 = 45041
$ws.Cells.Item(817, 5).Value2 = 7
$ws.Cells.Item(817, 6).Value2 = 100112004
$ws.Cells.Item(817, 7).Value2 = "Cebolla"
$ws.Cells.Item(817, 8).Value2 = "Sin especificar"
$ws.Cells.Item(817, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(817, 10).Value2 = 1500
$ws.Cells.Item(817, 11).Value2 = 7000
$ws.Cells.Item(817, 12).Value2 = 7000
$ws.Cells.Item(817, 13).Value2 = 7000
$ws.Cells.Item(817, 14).Value2 = "`$/malla 17 kilos"
$ws.Cells.Item(817, 15).Value2 = "Región del Maule"
$ws.Cells.Item(817, 16).Value2 = 412
$ws.Cells.Item(817, 17).Value2 = 17
$ws.Cells.Item(817, 18).Value2 = "Hortaliza"

# --- New row 818 ---
$ws.Cells.Item(818, 1).Value2 = 5
$ws.Cells.Item(818, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(818, 3).Value2 = "Maule"
$ws.Cells.Item(818, 4).Value2 = 45041
$ws.Cells.Item(818, 5).Value2 = 7
$ws.Cells.Item(818, 6).Value2 = 100112004
$ws.Cells.Item(818, 7).Value2 = "Cebolla"
$ws.Cells.Item(818, 8).Value2 = "Sin especificar"
$ws.Cells.Item(818, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(818, 10).Value2 = 1500
$ws.Cells.Item(818, 11).Value2 = 10000
$ws.Cells.Item(818, 12).Value2 = 10000
$ws.Cells.Item(818, 13).Value2 = 10000
$ws.Cells.Item(818, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(818, 15).Value2 = "Región del Maule"
$ws.Cells.Item(818, 16).Value2 = 400
$ws.Cells.Item(818, 17).Value2 = 25
$ws.Cells.Item(818, 18).Value2 = "Hortaliza"

# --- New row 819 ---
$ws.Cells.Item(819, 1).Value2 = 5
$ws.Cells.Item(819, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(819, 3).Value2 = "Maule"
$ws.Cells.Item(819, 4).Value2 = 45041
$ws.Cells.Item(819, 5).Value2 = 7
$ws.Cells.Item(819, 6).Value2 = 100112004
$ws.Cells.Item(819, 7).Value2 = "Cebolla"
$ws.Cells.Item(819, 8).Value2 = "Sin especificar"
$ws.Cells.Item(819, 9).Value2 = "2a (guarda)"
$ws.Cells.Item(819, 10).Value2 = 500
$ws.Cells.Item(819, 11).Value2 = 9000
$ws.Cells.Item(819, 12).Value2 = 9000
$ws.Cells.Item(819, 13).Value2 = 9000
$ws.Cells.Item(819, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(819, 15).Value2 = "Región del Maule"
$ws.Cells.Item(819, 16).Value2 = 360
$ws.Cells.Item(819, 17).Value2 = 25
$ws.Cells.Item(819, 18).Value2 = "Hortaliza"
